# Adds a new "2022-Q4" sheet (with its fund-holdings data) right after the
# "总计" summary sheet, shifting all the quarter sheets down by one
# position, and inserts a corresponding new row into the "总计" sheet's
# summary table.

function Set-TextCell($ws, $addr, $val) {
    # Forces a numeric-looking string to be stored as TEXT (inlineStr/shared
    # string) instead of being auto-coerced to a number by the .Value setter.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before the current "2022-Q3" tab.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($template)
$newSheet.Name = "2022-Q4"

# Seed the new sheet with the template's layout/styling (header row style,
# A-column style, dimensions) by copying its used range, then overwrite the
# cell contents with the real 2022-Q4 figures below.
$template.Range("A1:H13").Copy($newSheet.Range("A1:H13"))

# Clear out rows 8-13 (the template has 12 data rows, 2022-Q4 only has 6).
$newSheet.Range("A8:H13").Clear()

# -- header row (unchanged from template, already copied) --

# -- data rows --
$newSheet.Range("A2").Value = 0
Set-TextCell $newSheet "B2" "011429"
$newSheet.Range("C2").Value = "前海开源民裕进取混合"
Set-TextCell $newSheet "D2" "2.33"
Set-TextCell $newSheet "E2" "62.27"
Set-TextCell $newSheet "F2" "4.42"
Set-TextCell $newSheet "G2" "0.1030"
$newSheet.Range("H2").Value = 2

$newSheet.Range("A3").Value = 1
Set-TextCell $newSheet "B3" "630015"
$newSheet.Range("C3").Value = "华商大盘量化精选混合"
Set-TextCell $newSheet "D3" "2.48"
Set-TextCell $newSheet "E3" "89.21"
Set-TextCell $newSheet "F3" "4.00"
Set-TextCell $newSheet "G3" "0.0992"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
Set-TextCell $newSheet "B4" "008488"
$newSheet.Range("C4").Value = "华商恒益稳健混合"
Set-TextCell $newSheet "D4" "4.22"
Set-TextCell $newSheet "E4" "49.85"
Set-TextCell $newSheet "F4" "0.93"
Set-TextCell $newSheet "G4" "0.0392"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
Set-TextCell $newSheet "B5" "011588"
$newSheet.Range("C5").Value = "前海开源成份精选混合"
Set-TextCell $newSheet "D5" "0.77"
Set-TextCell $newSheet "E5" "62.53"
Set-TextCell $newSheet "F5" "4.20"
Set-TextCell $newSheet "G5" "0.0323"
$newSheet.Range("H5").Value = 3

$newSheet.Range("A6").Value = 4
Set-TextCell $newSheet "B6" "006890"
$newSheet.Range("C6").Value = "上投摩根领先优选混合A"
Set-TextCell $newSheet "D6" "0.29"
Set-TextCell $newSheet "E6" "81.41"
Set-TextCell $newSheet "F6" "3.39"
Set-TextCell $newSheet "G6" "0.0098"
$newSheet.Range("H6").Value = 5

$newSheet.Range("A7").Value = 5
Set-TextCell $newSheet "B7" "017098"
$newSheet.Range("C7").Value = "上投摩根领先优选混合C"
Set-TextCell $newSheet "D7" "0.00"
Set-TextCell $newSheet "E7" "81.41"
Set-TextCell $newSheet "F7" "3.39"
$newSheet.Range("G7").Value = 0
$newSheet.Range("H7").Value = 5

# Reset the "forced text" cells' formatting back to the sheet's default
# (no explicit style), matching the other plain data cells, by pasting in
# formats copied from a guaranteed-blank, unstyled cell.
$newSheet.Range("Z1").Copy()
$newSheet.Range("B2:G7").PasteSpecial(-4122)
$newSheet.Range("Z1").Clear()

# ---------------------------------------------------------------------
# 2. Insert the matching summary row into "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Seed the new row's formatting from the row right below it (which now
# holds the old row-2 data, shifted down), then overwrite with the new
# 2022-Q4 figures.
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2:D2"))
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.28
